# Natmi following Dr Hou advice
# Updates the NATMI LR-pair stats (Pltp-Abca1) on the active sheet: the
# "expressing cells" counts (Ligand: col E, Receptor: col K) change from 1 to 3,
# and all of the dependent expression / specificity / edge-weight metrics
# (G,H,I,J,M,N,O,P,Q,R,S,T) are refreshed with their recomputed values for rows 2-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns being updated, in order, for every data row (header row 1 is untouched):
#   E  Ligand-expressing cells
#   G  Ligand average expression value
#   H  Ligand total expression value
#   I  Ligand derived specificity of average expression value
#   J  Ligand derived specificity of total expression value
#   K  Receptor-expressing cells
#   M  Receptor average expression value
#   N  Receptor total expression value
#   O  Receptor derived specificity of average expression value
#   P  Receptor derived specificity of total expression value
#   Q  Edge average expression weight
#   R  Edge total expression weight
#   S  Edge average expression derived specificity
#   T  Edge total expression derived specificity
$columns = @(5, 7, 8, 9, 10, 11, 13, 14, 15, 16, 17, 18, 19, 20)

# One row of new values per worksheet row 2..17, in the same column order as above.
$newValues = @(
    @(3, 16.43128366666667, 49.293851, 0.1589154288049893, 0.1589154288049892, 3, 50.159191, 150.477573, 0.3080127083327998, 0.3080127083327999, 824.1798958115138, 7417.619062303624, 0.04894797162209297, 0.04894797162209297),
    @(3, 16.43128366666667, 49.293851, 0.1589154288049893, 0.1589154288049892, 3, 35.41899766666667, 106.256993, 0.2174975548896536, 0.2174975548896536, 581.9795978500049, 5237.816380650043, 0.03456371719932599, 0.03456371719932599),
    @(3, 16.43128366666667, 49.293851, 0.1589154288049893, 0.1589154288049892, 3, 66.72104766666666, 200.163143, 0.4097141557688168, 0.4097141557688168, 1096.312460748188, 9866.812146733693, 0.06510990075147569, 0.06510990075147569),
    @(3, 16.43128366666667, 49.293851, 0.1589154288049893, 0.1589154288049892, 3, 10.54856066666667, 31.645682, 0.06477558100872968, 0.06477558100872968, 173.3263925890425, 1559.937533301382, 0.0102938392320946, 0.01029383923209459),
    @(3, 9.522468333333334, 28.567405, 0.09209670827748465, 0.09209670827748463, 3, 50.159191, 150.477573, 0.3080127083327998, 0.3080127083327999, 477.6393079231184, 4298.753771308066, 0.02836695654508383, 0.02836695654508383),
    @(3, 9.522468333333334, 28.567405, 0.09209670827748465, 0.09209670827748463, 3, 35.41899766666667, 106.256993, 0.2174975548896536, 0.2174975548896536, 337.2762836792406, 3035.486553113165, 0.02003080886373863, 0.02003080886373863),
    @(3, 9.522468333333334, 28.567405, 0.09209670827748465, 0.09209670827748463, 3, 66.72104766666666, 200.163143, 0.4097141557688168, 0.4097141557688168, 635.3490635726572, 5718.141572153915, 0.03773332508099662, 0.03773332508099662),
    @(3, 9.522468333333334, 28.567405, 0.09209670827748465, 0.09209670827748463, 3, 10.54856066666667, 31.645682, 0.06477558100872968, 0.06477558100872968, 100.4483349105789, 904.03501419521, 0.005965617787665552, 0.005965617787665551),
    @(3, 76.05823733333334, 228.174712, 0.7355984867145992, 0.735598486714599, 3, 50.159191, 150.477573, 0.3080127083327998, 0.3080127083327999, 3815.019653525997, 34335.17688173398, 0.2265736821384728, 0.2265736821384728),
    @(3, 76.05823733333334, 228.174712, 0.7355984867145992, 0.735598486714599, 3, 35.41899766666667, 106.256993, 0.2174975548896536, 0.2174975548896536, 2693.906530640113, 24245.15877576102, 0.1599908722409547, 0.1599908722409546),
    @(3, 76.05823733333334, 228.174712, 0.7355984867145992, 0.735598486714599, 3, 66.72104766666666, 200.163143, 0.4097141557688168, 0.4097141557688168, 5074.685278559979, 45672.16750703981, 0.3013851129690912, 0.3013851129690912),
    @(3, 76.05823733333334, 228.174712, 0.7355984867145992, 0.735598486714599, 3, 10.54856066666667, 31.645682, 0.06477558100872968, 0.06477558100872968, 802.3049307103983, 7220.744376393584, 0.04764881936608048, 0.04764881936608047),
    @(3, 1.384413333333333, 4.15324, 0.01338937620292709, 0.01338937620292709, 3, 50.159191, 150.477573, 0.3080127083327998, 0.3080127083327999, 69.44105280961334, 624.96947528652, 0.004124098027150312, 0.004124098027150312),
    @(3, 1.384413333333333, 4.15324, 0.01338937620292709, 0.01338937620292709, 3, 35.41899766666667, 106.256993, 0.2174975548896536, 0.2174975548896536, 49.03453262303557, 441.3107936073201, 0.002912156585634356, 0.002912156585634356),
    @(3, 1.384413333333333, 4.15324, 0.01338937620292709, 0.01338937620292709, 3, 66.72104766666666, 200.163143, 0.4097141557688168, 0.4097141557688168, 92.36950800370222, 831.32557203332, 0.005485816967253358, 0.005485816967253357),
    @(3, 1.384413333333333, 4.15324, 0.01338937620292709, 0.01338937620292709, 3, 10.54856066666667, 31.645682, 0.06477558100872968, 0.06477558100872968, 14.60356803440889, 131.43211230968, 0.000867304622889061, 0.0008673046228890609)
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $newValues[$i]
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $ws.Cells.Item($row, $columns[$j]).Value = $rowValues[$j]
    }
}
